$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 54.69462833333333
$ws.Range("H2").Value = 164.083885
$ws.Range("I2").Value = 0.2790924419198448
$ws.Range("J2").Value = 0.2790924419198448
$ws.Range("M2").Value = 1.819857
$ws.Range("N2").Value = 5.459571
$ws.Range("O2").Value = 0.01485317462584607
$ws.Range("P2").Value = 0.01485317462584607
$ws.Range("Q2").Value = 99.536402234815
$ws.Range("R2").Value = 895.8276201133351
$ws.Range("S2").Value = 0.004145408776589255
$ws.Range("T2").Value = 0.004145408776589256
$ws.Range("G3").Value = 54.69462833333333
$ws.Range("H3").Value = 164.083885
$ws.Range("I3").Value = 0.2790924419198448
$ws.Range("J3").Value = 0.2790924419198448
$ws.Range("O3").Value = 0.726618572334523
$ws.Range("P3").Value = 0.7266185723345231
$ws.Range("Q3").Value = 4869.329305623533
$ws.Range("R3").Value = 43823.9637506118
$ws.Range("S3").Value = 0.2027937516971534
$ws.Range("T3").Value = 0.2027937516971534
$ws.Range("G4").Value = 54.69462833333333
$ws.Range("H4").Value = 164.083885
$ws.Range("I4").Value = 0.2790924419198448
$ws.Range("J4").Value = 0.2790924419198448
$ws.Range("M4").Value = 31.52924033333333
$ws.Range("N4").Value = 94.58772099999999
$ws.Range("O4").Value = 0.257333028084772
$ws.Range("P4").Value = 0.257333028084772
$ws.Range("Q4").Value = 1724.480081664009
$ws.Range("R4").Value = 15520.32073497608
$ws.Range("S4").Value = 0.07181970319480703
$ws.Range("T4").Value = 0.07181970319480703
$ws.Range("G5").Value = 54.69462833333333
$ws.Range("H5").Value = 164.083885
$ws.Range("I5").Value = 0.2790924419198448
$ws.Range("J5").Value = 0.2790924419198448
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.1464426666666667
$ws.Range("N5").Value = 0.439328
$ws.Range("O5").Value = 0.001195224954858853
$ws.Range("P5").Value = 0.001195224954858853
$ws.Range("Q5").Value = 8.009627225475555
$ws.Range("R5").Value = 72.08664502928001
$ws.Range("S5").Value = 0.0003335782512950934
$ws.Range("T5").Value = 0.0003335782512950935
$ws.Range("G6").Value = 19.32115333333334
$ws.Range("H6").Value = 57.96346000000001
$ws.Range("I6").Value = 0.09859081282432611
$ws.Range("J6").Value = 0.09859081282432611
$ws.Range("M6").Value = 1.819857
$ws.Range("N6").Value = 5.459571
$ws.Range("O6").Value = 0.01485317462584607
$ws.Range("P6").Value = 0.01485317462584607
$ws.Range("Q6").Value = 35.16173614174001
$ws.Range("R6").Value = 316.4556252756601
$ws.Range("S6").Value = 0.001464386559383819
$ws.Range("T6").Value = 0.00146438655938382
$ws.Range("G7").Value = 19.32115333333334
$ws.Range("H7").Value = 57.96346000000001
$ws.Range("I7").Value = 0.09859081282432611
$ws.Range("J7").Value = 0.09859081282432611
$ws.Range("O7").Value = 0.726618572334523
$ws.Range("P7").Value = 0.7266185723345231
$ws.Range("Q7").Value = 1720.115137652534
$ws.Range("S7").Value = 0.07163791565971202
$ws.Range("T7").Value = 0.07163791565971203
$ws.Range("G8").Value = 19.32115333333334
$ws.Range("H8").Value = 57.96346000000001
$ws.Range("I8").Value = 0.09859081282432611
$ws.Range("J8").Value = 0.09859081282432611
$ws.Range("M8").Value = 31.52924033333333
$ws.Range("N8").Value = 94.58772099999999
$ws.Range("O8").Value = 0.257333028084772
$ws.Range("P8").Value = 0.257333028084772
$ws.Range("Q8").Value = 609.1812869638512
$ws.Range("R8").Value = 5482.63158267466
$ws.Range("S8").Value = 0.02537067240542281
$ws.Range("T8").Value = 0.02537067240542281
$ws.Range("G9").Value = 19.32115333333334
$ws.Range("H9").Value = 57.96346000000001
$ws.Range("I9").Value = 0.09859081282432611
$ws.Range("J9").Value = 0.09859081282432611
$ws.Range("K9").Value = 1
$ws.Range("L9").Value = 0.3333333333333333
$ws.Range("M9").Value = 0.1464426666666667
$ws.Range("N9").Value = 0.439328
$ws.Range("O9").Value = 0.001195224954858853
$ws.Range("P9").Value = 0.001195224954858853
$ws.Range("Q9").Value = 2.82944121720889
$ws.Range("R9").Value = 25.46497095488001
$ws.Range("S9").Value = 0.0001178381998074527
$ws.Range("T9").Value = 0.0001178381998074528
$ws.Range("G10").Value = 11.023718
$ws.Range("H10").Value = 33.071154
$ws.Range("I10").Value = 0.05625116157486912
$ws.Range("J10").Value = 0.05625116157486911
$ws.Range("M10").Value = 1.819857
$ws.Range("N10").Value = 5.459571
$ws.Range("O10").Value = 0.01485317462584607
$ws.Range("P10").Value = 0.01485317462584607
$ws.Range("Q10").Value = 20.061590368326
$ws.Range("R10").Value = 180.554313314934
$ws.Range("S10").Value = 0.0008355083257782132
$ws.Range("T10").Value = 0.0008355083257782132
$ws.Range("G11").Value = 11.023718
$ws.Range("H11").Value = 33.071154
$ws.Range("I11").Value = 0.05625116157486912
$ws.Range("J11").Value = 0.05625116157486911
$ws.Range("O11").Value = 0.726618572334523
$ws.Range("P11").Value = 0.7266185723345231
$ws.Range("Q11").Value = 981.41471566808
$ws.Range("R11").Value = 8832.73244101272
$ws.Range("S11").Value = 0.04087313871568998
$ws.Range("T11").Value = 0.04087313871568998
$ws.Range("G12").Value = 11.023718
$ws.Range("H12").Value = 33.071154
$ws.Range("I12").Value = 0.05625116157486912
$ws.Range("J12").Value = 0.05625116157486911
$ws.Range("M12").Value = 31.52924033333333
$ws.Range("N12").Value = 94.58772099999999
$ws.Range("O12").Value = 0.257333028084772
$ws.Range("P12").Value = 0.257333028084772
$ws.Range("Q12").Value = 347.5694541888927
$ws.Range("R12").Value = 3128.125087700033
$ws.Range("S12").Value = 0.01447528174134684
$ws.Range("T12").Value = 0.01447528174134684
$ws.Range("G13").Value = 11.023718
$ws.Range("H13").Value = 33.071154
$ws.Range("I13").Value = 0.05625116157486912
$ws.Range("J13").Value = 0.05625116157486911
$ws.Range("K13").Value = 1
$ws.Range("L13").Value = 0.3333333333333333
$ws.Range("M13").Value = 0.1464426666666667
$ws.Range("N13").Value = 0.439328
$ws.Range("O13").Value = 0.001195224954858853
$ws.Range("P13").Value = 0.001195224954858853
$ws.Range("Q13").Value = 1.614342660501333
$ws.Range("R13").Value = 14.529083944512
$ws.Range("S13").Value = [double]"6.723279205408096E-05"
$ws.Range("T13").Value = [double]"6.723279205408097E-05"
$ws.Range("G14").Value = 110.9336623333333
$ws.Range("H14").Value = 332.800987
$ws.Range("I14").Value = 0.5660655836809599
$ws.Range("J14").Value = 0.5660655836809599
$ws.Range("M14").Value = 1.819857
$ws.Range("N14").Value = 5.459571
$ws.Range("O14").Value = 0.01485317462584607
$ws.Range("P14").Value = 0.01485317462584607
$ws.Range("Q14").Value = 201.883401932953
$ws.Range("R14").Value = 1816.950617396577
$ws.Range("S14").Value = 0.008407870964094778
$ws.Range("T14").Value = 0.008407870964094778
$ws.Range("G15").Value = 110.9336623333333
$ws.Range("H15").Value = 332.800987
$ws.Range("I15").Value = 0.5660655836809599
$ws.Range("J15").Value = 0.5660655836809599
$ws.Range("O15").Value = 0.726618572334523
$ws.Range("P15").Value = 0.7266185723345231
$ws.Range("Q15").Value = 9876.153279400571
$ws.Range("R15").Value = 88885.37951460514
$ws.Range("S15").Value = 0.4113137662619676
$ws.Range("T15").Value = 0.4113137662619676
$ws.Range("G16").Value = 110.9336623333333
$ws.Range("H16").Value = 332.800987
$ws.Range("I16").Value = 0.5660655836809599
$ws.Range("J16").Value = 0.5660655836809599
$ws.Range("M16").Value = 31.52924033333333
$ws.Range("N16").Value = 94.58772099999999
$ws.Range("O16").Value = 0.257333028084772
$ws.Range("P16").Value = 0.257333028084772
$ws.Range("Q16").Value = 3497.654100764513
$ws.Range("R16").Value = 31478.88690688062
$ws.Range("S16").Value = 0.1456673707431953
$ws.Range("T16").Value = 0.1456673707431953
$ws.Range("G17").Value = 110.9336623333333
$ws.Range("H17").Value = 332.800987
$ws.Range("I17").Value = 0.5660655836809599
$ws.Range("J17").Value = 0.5660655836809599
$ws.Range("K17").Value = 1
$ws.Range("L17").Value = 0.3333333333333333
$ws.Range("M17").Value = 0.1464426666666667
$ws.Range("N17").Value = 0.439328
$ws.Range("O17").Value = 0.001195224954858853
$ws.Range("P17").Value = 0.001195224954858853
$ws.Range("Q17").Value = 16.24542133519289
$ws.Range("R17").Value = 146.208792016736
$ws.Range("S17").Value = 0.0006765757117022254
$ws.Range("T17").Value = 0.0006765757117022255
